$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Sheets("ALC")
$ws.Range("H39").Value = 766.64703
$ws.Range("I39").Value = 855.9286
$ws.Range("J39").Value = 350
$ws.Range("K39").Value = 2567.7858
$ws.Range("L39").Value = 1050
$ws.Range("M39").Value = -2271.7858
$ws.Range("N39").Value = -1642
$ws.Range("H40").Value = 2168.875
$ws.Range("I40").Value = 2030
$ws.Range("J40").Value = 2400.3333
$ws.Range("K40").Value = 2030
$ws.Range("L40").Value = 2400.3333
$ws.Range("M40").Value = -1855
$ws.Range("N40").Value = -2750.3333
$ws.Range("H62").Value = 8435.727999999999
$ws.Range("I62").Value = 7332.5
$ws.Range("J62").Value = 9759.6
$ws.Range("K62").Value = 7332.5
$ws.Range("L62").Value = 9759.6
$ws.Range("M62").Value = -6708.5
$ws.Range("N62").Value = -11007.6
$ws.Range("H65").Value = 8435.727999999999
$ws.Range("I65").Value = 7332.5
$ws.Range("J65").Value = 9759.6
$ws.Range("K65").Value = 36662.5
$ws.Range("L65").Value = 48798
$ws.Range("M65").Value = -33542.5
$ws.Range("N65").Value = -55038
$ws.Range("H123").Value = 97099.625
$ws.Range("J123").Value = 97099.625
$ws.Range("L123").Value = 97099.625
$ws.Range("N123").Value = -106899.625
$ws.Range("H141").Value = 1794.129
$ws.Range("I141").Value = 1526.5927
$ws.Range("K141").Value = 4579.7781
$ws.Range("M141").Value = 600.2219000000005

# ---- Sheet: ARM ----
$ws = $wb.Sheets("ARM")
$ws.Range("H32").Value = 2480.6025
$ws.Range("I32").Value = 1756.5593
$ws.Range("J32").Value = 4728.9473
$ws.Range("K32").Value = 1756.5593
$ws.Range("L32").Value = 4728.9473
$ws.Range("M32").Value = -1469.5593
$ws.Range("N32").Value = -5302.9473
$ws.Range("H45").Value = 1971.5555
$ws.Range("I45").Value = 1422.2307
$ws.Range("J45").Value = 3399.8
$ws.Range("K45").Value = 1422.2307
$ws.Range("L45").Value = 3399.8
$ws.Range("M45").Value = -1045.2307
$ws.Range("N45").Value = -4153.8
$ws.Range("H74").Value = 1121.925
$ws.Range("I74").Value = 1229.9584
$ws.Range("J74").Value = 959.875
$ws.Range("K74").Value = 1229.9584
$ws.Range("L74").Value = 959.875
$ws.Range("M74").Value = -355.9584
$ws.Range("N74").Value = -2707.875
$ws.Range("H77").Value = 1121.925
$ws.Range("I77").Value = 1229.9584
$ws.Range("J77").Value = 959.875
$ws.Range("K77").Value = 6149.791999999999
$ws.Range("L77").Value = 4799.375
$ws.Range("M77").Value = -1781.791999999999
$ws.Range("N77").Value = -13535.375
$ws.Range("H132").Value = 2219.5789
$ws.Range("I132").Value = 1698.125
$ws.Range("J132").Value = 5000.6665
$ws.Range("K132").Value = 5094.375
$ws.Range("L132").Value = 15001.9995
$ws.Range("M132").Value = -2564.375
$ws.Range("N132").Value = -20061.9995

# ---- Sheet: BSM ----
$ws = $wb.Sheets("BSM")
$ws.Range("H105").Value = 17547046
$ws.Range("I105").Value = 20836454
$ws.Range("J105").Value = 3533.3333
$ws.Range("K105").Value = 20836454
$ws.Range("L105").Value = 3533.3333
$ws.Range("M105").Value = -20834707
$ws.Range("N105").Value = -7027.3333

# ---- Sheet: CRP ----
$ws = $wb.Sheets("CRP")
$ws.Range("H31").Value = 1322.7979
$ws.Range("I31").Value = 1009.125
$ws.Range("J31").Value = 1555.1482
$ws.Range("K31").Value = 1009.125
$ws.Range("L31").Value = 1555.1482
$ws.Range("M31").Value = -714.125
$ws.Range("N31").Value = -2145.1482
$ws.Range("H34").Value = 1322.7979
$ws.Range("I34").Value = 1009.125
$ws.Range("J34").Value = 1555.1482
$ws.Range("K34").Value = 1009.125
$ws.Range("L34").Value = 1555.1482
$ws.Range("M34").Value = -807.125
$ws.Range("N34").Value = -1959.1482
$ws.Range("H99").Value = 12501460
$ws.Range("I99").Value = 15626325
$ws.Range("J99").Value = 2000
$ws.Range("K99").Value = 15626325
$ws.Range("L99").Value = 2000
$ws.Range("M99").Value = -15624827
$ws.Range("N99").Value = -4996
$ws.Range("H105").Value = 426.23077
$ws.Range("I105").Value = 316.1
$ws.Range("K105").Value = 316.1
$ws.Range("M105").Value = 1430.9
$ws.Range("H122").Value = 1066.6666
$ws.Range("I122").Value = 800
$ws.Range("K122").Value = 2400
$ws.Range("M122").Value = 50
$ws.Range("H126").Value = 12501460
$ws.Range("I126").Value = 15626325
$ws.Range("J126").Value = 2000
$ws.Range("K126").Value = 46878975
$ws.Range("L126").Value = 6000
$ws.Range("M126").Value = -46876505
$ws.Range("N126").Value = -10940
$ws.Range("H132").Value = 2093.652
$ws.Range("I132").Value = 1529.625
$ws.Range("J132").Value = 3382.8572
$ws.Range("K132").Value = 4588.875
$ws.Range("L132").Value = 10148.5716
$ws.Range("M132").Value = -2058.875
$ws.Range("N132").Value = -15208.5716
$ws.Range("H134").Value = 2075.617
$ws.Range("I134").Value = 929
$ws.Range("J134").Value = 4778.357
$ws.Range("K134").Value = 2787
$ws.Range("L134").Value = 14335.071
$ws.Range("M134").Value = -252
$ws.Range("N134").Value = -19405.071

# ---- Sheet: CUL ----
$ws = $wb.Sheets("CUL")
$ws.Range("H39").Value = 8748.471
$ws.Range("J39").Value = 8748.471
$ws.Range("L39").Value = 26245.413
$ws.Range("N39").Value = -26833.413
$ws.Range("H68").Value = 827.23956
$ws.Range("J68").Value = 955.2727
$ws.Range("L68").Value = 2865.8181
$ws.Range("N68").Value = -4487.8181
$ws.Range("H71").Value = 827.23956
$ws.Range("J71").Value = 955.2727
$ws.Range("L71").Value = 8597.454299999999
$ws.Range("N71").Value = -16709.4543
$ws.Range("H81").Value = 5179776
$ws.Range("J81").Value = 9495444
$ws.Range("L81").Value = 28486332
$ws.Range("N81").Value = -28488578
$ws.Range("H84").Value = 5179776
$ws.Range("J84").Value = 9495444
$ws.Range("L84").Value = 85458996
$ws.Range("N84").Value = -85470228

# ---- Sheet: GSM ----
$ws = $wb.Sheets("GSM")
$ws.Range("H80").Value = 2000
$ws.Range("I80").Value = 1733.3334
$ws.Range("J80").Value = 2400
$ws.Range("K80").Value = 1733.3334
$ws.Range("L80").Value = 2400
$ws.Range("M80").Value = -735.3334
$ws.Range("N80").Value = -4396
$ws.Range("H83").Value = 2000
$ws.Range("I83").Value = 1733.3334
$ws.Range("J83").Value = 2400
$ws.Range("K83").Value = 8666.666999999999
$ws.Range("L83").Value = 12000
$ws.Range("M83").Value = -3674.666999999999
$ws.Range("N83").Value = -21984
$ws.Range("H102").Value = 2551.647
$ws.Range("I102").Value = 1460.3636
$ws.Range("J102").Value = 4552.3335
$ws.Range("K102").Value = 1460.3636
$ws.Range("L102").Value = 4552.3335
$ws.Range("M102").Value = 161.6364000000001
$ws.Range("N102").Value = -7796.3335
$ws.Range("H122").Value = 1389915.4
$ws.Range("I122").Value = 1852620.5
$ws.Range("J122").Value = 1800
$ws.Range("K122").Value = 5557861.5
$ws.Range("L122").Value = 5400
$ws.Range("M122").Value = -5555411.5
$ws.Range("N122").Value = -10300
$ws.Range("H126").Value = 1977.7354
$ws.Range("I126").Value = 1683.6666
$ws.Range("J126").Value = 2138.1365
$ws.Range("K126").Value = 5050.9998
$ws.Range("L126").Value = 6414.4095
$ws.Range("M126").Value = -2580.9998
$ws.Range("N126").Value = -11354.4095

# ---- Sheet: LTW ----
$ws = $wb.Sheets("LTW")
$ws.Range("H87").Value = 500189
$ws.Range("J87").Value = 500189
$ws.Range("L87").Value = 500189
$ws.Range("N87").Value = -502435
$ws.Range("H90").Value = 500189
$ws.Range("J90").Value = 500189
$ws.Range("L90").Value = 1500567
$ws.Range("N90").Value = -1511799
$ws.Range("H100").Value = 2346.5356
$ws.Range("I100").Value = 1625.9
$ws.Range("K100").Value = 1625.9
$ws.Range("M100").Value = -1084.9
$ws.Range("H122").Value = 3465.36
$ws.Range("I122").Value = 1490
$ws.Range("J122").Value = 3959.2
$ws.Range("K122").Value = 4470
$ws.Range("L122").Value = 11877.6
$ws.Range("M122").Value = -2020
$ws.Range("N122").Value = -16777.6
$ws.Range("H132").Value = 2864.2856
$ws.Range("I132").Value = 2227.24
$ws.Range("J132").Value = 3801.1177
$ws.Range("K132").Value = 6681.719999999999
$ws.Range("L132").Value = 11403.3531
$ws.Range("M132").Value = -4151.719999999999
$ws.Range("N132").Value = -16463.3531
$ws.Range("H135").Value = 40427.832
$ws.Range("J135").Value = 40427.832
$ws.Range("L135").Value = 40427.832
$ws.Range("N135").Value = -50567.832
$ws.Range("H136").Value = 5459.04
$ws.Range("I136").Value = 1978.2667
$ws.Range("J136").Value = 10680.2
$ws.Range("K136").Value = 5934.800099999999
$ws.Range("L136").Value = 32040.6
$ws.Range("M136").Value = -3384.800099999999
$ws.Range("N136").Value = -37140.60000000001

# ---- Sheet: WVR ----
$ws = $wb.Sheets("WVR")
$ws.Range("H100").Value = 300
$ws.Range("I100").Value = 300
$ws.Range("J100").Value = 0
$ws.Range("K100").Value = 600
$ws.Range("L100").Value = 0
$ws.Range("M100").Value = -59
$ws.Range("N100").ClearContents()
$ws.Range("H126").Value = 91947.82000000001
$ws.Range("I126").Value = 112094.555
$ws.Range("J126").Value = 1287.5
$ws.Range("K126").Value = 336283.665
$ws.Range("L126").Value = 3862.5
$ws.Range("M126").Value = -333813.665
$ws.Range("N126").Value = -8802.5
$ws.Range("H132").Value = 11906424
$ws.Range("I132").Value = 14707125
$ws.Range("J132").Value = 3443.125
$ws.Range("K132").Value = 44121375
$ws.Range("L132").Value = 10329.375
$ws.Range("M132").Value = -44118845
$ws.Range("N132").Value = -15389.375
$ws.Range("H136").Value = 17598128
$ws.Range("I136").Value = 30394638
$ws.Range("J136").Value = 2927.875
$ws.Range("K136").Value = 91183914
$ws.Range("L136").Value = 8783.625
$ws.Range("M136").Value = -91181364
$ws.Range("N136").Value = -13883.625
Write-Host "Edit complete"
